$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 2014
$ws.Range("I2").Value = 2998.5
$ws.Range("J2").Value = 701.3333
$ws.Range("K2").Value = 2998.5
$ws.Range("L2").Value = 701.3333
$ws.Range("M2").Value = -2885.5
$ws.Range("N2").Value = -927.3333
$ws.Range("H18").Value = 8462.073
$ws.Range("I18").Value = 8031.567
$ws.Range("K18").Value = 8031.567
$ws.Range("M18").Value = -7747.567
$ws.Range("H111").Value = 1199.8572
$ws.Range("I111").Value = 579.8
$ws.Range("K111").Value = 1739.4
$ws.Range("M111").Value = 1327.6
$ws.Range("H113").Value = 2753.3635
$ws.Range("I113").Value = 2665.6667
$ws.Range("J113").Value = 2858.6
$ws.Range("K113").Value = 2665.6667
$ws.Range("L113").Value = 2858.6
$ws.Range("M113").Value = 588.3332999999998
$ws.Range("N113").Value = -9366.6
$ws.Range("H116").Value = 7333.3335
$ws.Range("I116").Value = 7500
$ws.Range("K116").Value = 7500
$ws.Range("M116").Value = -4058
$ws.Range("H118").Value = 89.875
$ws.Range("I118").Value = 95.57143000000001
$ws.Range("J118").Value = 50
$ws.Range("K118").Value = 286.71429
$ws.Range("L118").Value = 150
$ws.Range("M118").Value = 1370.28571
$ws.Range("N118").Value = -3464
$ws.Range("H138").Value = 3610.3076
$ws.Range("I138").Value = 4799
$ws.Range("J138").Value = 2981
$ws.Range("K138").Value = 14397
$ws.Range("L138").Value = 8943
$ws.Range("M138").Value = -9257
$ws.Range("N138").Value = -19223

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H18").Value = 0
$ws.Range("J18").Value = 0
$ws.Range("L18").Value = 0
$ws.Range("N18").ClearContents()
$ws.Range("H32").Value = 235955.98
$ws.Range("I32").Value = 2411.054
$ws.Range("K32").Value = 2411.054
$ws.Range("M32").Value = -2124.054
$ws.Range("H110").Value = 1305.2142
$ws.Range("I110").Value = 636.5333000000001
$ws.Range("J110").Value = 2076.7693
$ws.Range("K110").Value = 636.5333000000001
$ws.Range("L110").Value = 2076.7693
$ws.Range("M110").Value = 1408.4667
$ws.Range("N110").Value = -6166.7693
$ws.Range("H132").Value = 2999.5
$ws.Range("J132").Value = 2999
$ws.Range("L132").Value = 8997
$ws.Range("N132").Value = -14057

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 914878.4
$ws.Range("I20").Value = 7150.8335
$ws.Range("K20").Value = 7150.8335
$ws.Range("M20").Value = -6903.8335
$ws.Range("H107").Value = 1661.0555
$ws.Range("I107").Value = 1075.4286
$ws.Range("K107").Value = 1075.4286
$ws.Range("M107").Value = 844.5714
$ws.Range("H134").Value = 2434.3333
$ws.Range("I134").Value = 2174.2307
$ws.Range("K134").Value = 6522.6921
$ws.Range("M134").Value = -3987.6921

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1092.3
$ws.Range("I16").Value = 1069.2
$ws.Range("K16").Value = 1069.2
$ws.Range("M16").Value = -782.2
$ws.Range("H31").Value = 11368
$ws.Range("I31").Value = 15097.286
$ws.Range("K31").Value = 15097.286
$ws.Range("M31").Value = -14802.286
$ws.Range("H34").Value = 11368
$ws.Range("I34").Value = 15097.286
$ws.Range("K34").Value = 15097.286
$ws.Range("M34").Value = -14895.286
$ws.Range("H74").Value = 34279.8
$ws.Range("J74").Value = 34279.8
$ws.Range("L74").Value = 34279.8
$ws.Range("N74").Value = -36027.8
$ws.Range("H77").Value = 34279.8
$ws.Range("J77").Value = 34279.8
$ws.Range("L77").Value = 102839.4
$ws.Range("N77").Value = -111575.4
$ws.Range("H99").Value = 4770.625
$ws.Range("I99").Value = 4770.625
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 4770.625
$ws.Range("L99").Value = 0
$ws.Range("M99").Value = -3272.625
$ws.Range("N99").ClearContents()
$ws.Range("H105").Value = 721.3333
$ws.Range("I105").Value = 974.2857
$ws.Range("K105").Value = 974.2857
$ws.Range("M105").Value = 772.7143
$ws.Range("H106").Value = 25990
$ws.Range("J106").Value = 25990
$ws.Range("L106").Value = 25990
$ws.Range("N106").Value = -28514
$ws.Range("H107").Value = 1150.3438
$ws.Range("I107").Value = 841.9
$ws.Range("J107").Value = 1290.5454
$ws.Range("K107").Value = 841.9
$ws.Range("L107").Value = 1290.5454
$ws.Range("M107").Value = 1078.1
$ws.Range("N107").Value = -5130.5454
$ws.Range("H113").Value = 1092.3
$ws.Range("I113").Value = 1069.2
$ws.Range("K113").Value = 1069.2
$ws.Range("M113").Value = 1100.8
$ws.Range("H126").Value = 4770.625
$ws.Range("I126").Value = 4770.625
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 14311.875
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -11841.875
$ws.Range("N126").ClearContents()
$ws.Range("H132").Value = 5814.6875
$ws.Range("I132").Value = 4462.4287
$ws.Range("K132").Value = 13387.2861
$ws.Range("M132").Value = -10857.2861
$ws.Range("H134").Value = 2210.75
$ws.Range("I134").Value = 1503
$ws.Range("K134").Value = 4509
$ws.Range("M134").Value = -1974

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()
$ws.Range("N5").ClearContents()
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("M124").ClearContents()
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 37755.1
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()
$ws.Range("H24").Value = 14217.097
$ws.Range("I24").Value = 4967.684
$ws.Range("K24").Value = 4967.684
$ws.Range("M24").Value = -4794.684
$ws.Range("H132").Value = 5171.1665
$ws.Range("I132").Value = 4999
$ws.Range("J132").Value = 5205.6
$ws.Range("K132").Value = 14997
$ws.Range("L132").Value = 15616.8
$ws.Range("M132").Value = -12467
$ws.Range("N132").Value = -20676.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1750
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 1750
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 1750
$ws.Range("M22").ClearContents()
$ws.Range("N22").Value = -2340
$ws.Range("H27").Value = 1750
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 1750
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 1750
$ws.Range("M27").ClearContents()
$ws.Range("N27").Value = -1964
$ws.Range("H55").Value = 743.86664
$ws.Range("I55").Value = 743.7692
$ws.Range("K55").Value = 743.7692
$ws.Range("M55").Value = -570.7692
$ws.Range("H82").Value = 101204.9
$ws.Range("I82").Value = 1528.1428
$ws.Range("J82").Value = 333784
$ws.Range("K82").Value = 1528.1428
$ws.Range("L82").Value = 333784
$ws.Range("M82").Value = -1167.1428
$ws.Range("N82").Value = -334506
$ws.Range("H85").Value = 101204.9
$ws.Range("I85").Value = 1528.1428
$ws.Range("J85").Value = 333784
$ws.Range("K85").Value = 1528.1428
$ws.Range("L85").Value = 333784
$ws.Range("M85").Value = -280.1428000000001
$ws.Range("N85").Value = -336280
$ws.Range("H132").Value = 260463.25
$ws.Range("I132").Value = 333951
$ws.Range("K132").Value = 1001853
$ws.Range("M132").Value = -999323

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 7000
$ws.Range("J39").Value = 7000
$ws.Range("L39").Value = 7000
$ws.Range("N39").Value = -7826
$ws.Range("H62").Value = 7036.6665
$ws.Range("I62").Value = 5960
$ws.Range("K62").Value = 5960
$ws.Range("M62").Value = -5336
$ws.Range("H65").Value = 7036.6665
$ws.Range("I65").Value = 5960
$ws.Range("K65").Value = 29800
$ws.Range("M65").Value = -26680
$ws.Range("H104").Value = 17711.625
$ws.Range("J104").Value = 17711.625
$ws.Range("L104").Value = 17711.625
$ws.Range("N104").Value = -24699.625
$ws.Range("H113").Value = 380.33334
$ws.Range("I113").Value = 281.92856
$ws.Range("J113").Value = 724.75
$ws.Range("K113").Value = 845.78568
$ws.Range("L113").Value = 2174.25
$ws.Range("M113").Value = 1324.21432
$ws.Range("N113").Value = -6514.25
